# Append a new row (44) of data to each of the 4 worksheets, duplicating the
# last existing row (43) but with an updated timestamp in column A.
#
# Columns: A=time, B=总长, C=ID, D=实际长度, E=和校验,
#          F=总长_DEC, G=ID_DEC, H=实际长度_DEC, I=和校验_DEC

$wb = $excel.ActiveWorkbook

$newTime = 45830.46016203704

$rows = @(
    @{
        Sheet = "MID_LFT_#1"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x74"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 372
        I = 7
    },
    @{
        Sheet = "MID_LFT_#2"
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x68"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 360
        I = 25
    },
    @{
        Sheet = "MID_PLT_#1"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x6A"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 106
        I = 15
    },
    @{
        Sheet = "MID_PLT_#2"
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7F"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 127
        I = 9
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $newRowNum = 44

    $ws.Cells.Item($newRowNum, 1).Value = $newTime
    $ws.Cells.Item($newRowNum, 1).NumberFormat = $ws.Cells.Item($newRowNum - 1, 1).NumberFormat

    $ws.Cells.Item($newRowNum, 2).Value = $row.B
    $ws.Cells.Item($newRowNum, 3).Value = $row.C
    $ws.Cells.Item($newRowNum, 4).Value = $row.D
    $ws.Cells.Item($newRowNum, 5).Value = $row.E
    $ws.Cells.Item($newRowNum, 6).Value = $row.F
    $ws.Cells.Item($newRowNum, 7).Value = $row.G
    $ws.Cells.Item($newRowNum, 8).Value = $row.H
    $ws.Cells.Item($newRowNum, 9).Value = $row.I
}
